# Updates Sheet1!D2:E51 (Price / Volume(1h) columns) with refreshed
# cryptocurrency quote data, matching the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.593.64"
$ws.Range("D3").Value = "1.667.49"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'236.98"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4803"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.2631"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "'0.06163"
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("D10").Value = "'0.07096"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").Value = "1.666.26"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "'14.90"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("D13").Value = "'0.5989"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "'74.72"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "25.583.18"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").Value = "'0.000006802"
$ws.Range("E19").Value = "  +4.40%  "
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "'4.478"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "1.878.10"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "'8.713"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").Value = "'5.368"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("D25").Value = "'134.34"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "'15.12"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'104.86"
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'3.984"
$ws.Range("E30").Value = "  +4.40%  "
$ws.Range("D31").Value = "'3.675"
$ws.Range("E31").Value = "  +4.40%  "
$ws.Range("D32").Value = "'0.07702"
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("D33").Value = "'0.04361"
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("D34").Value = "'0.9998"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'2.617"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").Value = "'0.6163"
$ws.Range("E36").Value = "  +5.61%  "
$ws.Range("D37").Value = "'0.9528"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").Value = "'2.608"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "'0.8689"
$ws.Range("E39").Value = "  +3.58%  "
$ws.Range("D40").Value = "'1.000"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "'0.01520"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "'1.874"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("D43").Value = "'97.86"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "'0.3783"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("D45").Value = "'4.691"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").Value = "'0.1127"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'6.247"
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("D49").Value = "'29.61"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").Value = "'7.430"
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("D51").Value = "'0.3357"
$ws.Range("E51").Value = "  +1.42%  "
